# Actualización automática 2025-09-08 09:55:08
#
# Insert a new salesperson "GARCIA BRAVO JOSE LUIS" (all-zero row) at row 14
# in the two detail sheets ("VENTAS POR GRUPO" and "VENTA MENSUAL"),
# pushing the existing rows 14-26 down to 15-27 and the totals row from
# 27 -> 28. Then refresh the "X de 25" -> "X de 26" counters on the new
# totals row, and update the CUMPLIMIENTO MENSUAL summary sheet (values +
# two column widths).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:R, data rows 2-26, totals 27)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(14).Insert()
$ws1.Cells.Item(14, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(14, 2).Value = "GARCIA BRAVO JOSE LUIS"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(14, $c).Value = 0
}

# refresh the "X de 25" -> "X de 26" counters on the (now shifted) totals row 28
$ws1.Cells.Item(28, 3).Value  = "0 de 26"
$ws1.Cells.Item(28, 4).Value  = "1 de 26"
$ws1.Cells.Item(28, 5).Value  = "1 de 26"
$ws1.Cells.Item(28, 6).Value  = "0 de 26"
$ws1.Cells.Item(28, 7).Value  = "0 de 26"
$ws1.Cells.Item(28, 8).Value  = "0 de 26"
$ws1.Cells.Item(28, 9).Value  = "0 de 26"
$ws1.Cells.Item(28, 10).Value = "0 de 26"
$ws1.Cells.Item(28, 11).Value = "0 de 26"
$ws1.Cells.Item(28, 12).Value = "1 de 26"
$ws1.Cells.Item(28, 13).Value = "3 de 26"
$ws1.Cells.Item(28, 14).Value = "0 de 26"
$ws1.Cells.Item(28, 15).Value = "0 de 26"
$ws1.Cells.Item(28, 16).Value = "0 de 26"
$ws1.Cells.Item(28, 17).Value = "0 de 26"
$ws1.Cells.Item(28, 18).Value = "0 de 26"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:G, data rows 2-26, totals 27)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(14).Insert()
$ws2.Cells.Item(14, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(14, 2).Value = "GARCIA BRAVO JOSE LUIS"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(14, $c).Value = 0
}

# ---------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL" -- refreshed summary figures + widths
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Cells.Item(2, 4).Value = 4083.42
$ws3.Cells.Item(2, 5).Value = -4083.42

$ws3.Cells.Item(3, 3).Value = 20000
$ws3.Cells.Item(3, 4).Value = 9694.08
$ws3.Cells.Item(3, 5).Value = 10305.92
$ws3.Cells.Item(3, 6).Value = 0.484704

$ws3.Cells.Item(4, 3).Value = 20000
$ws3.Cells.Item(4, 4).Value = 13777.5
$ws3.Cells.Item(4, 5).Value = 6222.5
$ws3.Cells.Item(4, 6).Value = 0.688875

# column D: raw OOXML width 14 -> 13 ; column F: raw OOXML width 24 -> 18
# (the ColumnWidth COM property is offset from the stored <col width> by
# ~0.8333 character units, so subtract that back out)
$ws3.Columns.Item(4).ColumnWidth = 13 - 0.8333333333333333
$ws3.Columns.Item(6).ColumnWidth = 18 - 0.8333333333333333
